$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = $origStyle
}

$ws.Range("D2").Value = '42.165.04'
$ws.Range("E2").Value = '  +1.30%  '

$ws.Range("D3").Value = '2.165.93'
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("E4").Value = '  -0.04%  '

Set-TextValue $ws.Range("D5") '253.04'
$ws.Range("E5").Value = '  +6.30%  '

$ws.Range("E6").Value = '  +0.13%  '

Set-TextValue $ws.Range("D7") '73.05'
$ws.Range("E7").Value = '  +1.47%  '

$ws.Range("E8").Value = '  +0.02%  '

Set-TextValue $ws.Range("D9") '0.574'
$ws.Range("E9").Value = '  -0.68%  '

Set-TextValue $ws.Range("D10") '39.59'
$ws.Range("E10").Value = '  -0.30%  '

Set-TextValue $ws.Range("D11") '0.0906'
$ws.Range("E11").Value = '  -0.10%  '

$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("D14").Value = '2.488.69'
$ws.Range("E14").Value = '  +0.03%  '

Set-TextValue $ws.Range("D15") '14.13'
$ws.Range("E15").Value = '  -1.04%  '

$ws.Range("D16").Value = '2.172.94'
$ws.Range("E16").Value = '  +1.48%  '

Set-TextValue $ws.Range("D17") '0.761'
$ws.Range("E17").Value = '  -2.00%  '

$ws.Range("D18").Value = '42.060.35'
$ws.Range("E18").Value = '  +1.41%  '

$ws.Range("E19").Value = '  -0.68%  '

Set-TextValue $ws.Range("D20") '70.34'
$ws.Range("E20").Value = '  +0.52%  '

Set-TextValue $ws.Range("D21") '5.80'
$ws.Range("E21").Value = '  +0.44%  '

$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D22") '9.54'
$ws.Range("E22").Value = '  -3.79%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D23") '225.74'
$ws.Range("E23").Value = '  -0.12%  '

Set-TextValue $ws.Range("D24") '2.13'
$ws.Range("E24").Value = '  +6.53%  '

$ws.Range("E25").Value = '  -0.12%  '

Set-TextValue $ws.Range("D26") '10.41'
$ws.Range("E26").Value = '  -2.96%  '

$ws.Range("E27").Value = '  +1.32%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D28") '2.22'
$ws.Range("E28").Value = '  +2.72%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D29") '2.18'
$ws.Range("E29").Value = '  -0.32%  '

Set-TextValue $ws.Range("D30") '36.38'
$ws.Range("E30").Value = '  +10.84%  '

Set-TextValue $ws.Range("D31") '167.99'
$ws.Range("E31").Value = '  -1.84%  '

Set-TextValue $ws.Range("D32") '19.90'
$ws.Range("E32").Value = '  +0.51%  '

Set-TextValue $ws.Range("D33") '0.0804'
$ws.Range("E33").Value = '  +4.32%  '

$ws.Range("E34").Value = '  -5.17%  '

$ws.Range("E35").Value = '  -0.57%  '

Set-TextValue $ws.Range("D36") '0.106'
$ws.Range("E36").Value = '  +4.01%  '

Set-TextValue $ws.Range("D37") '4.22'
$ws.Range("E37").Value = '  -1.39%  '

$ws.Range("E38").Value = '  +8.86%  '

Set-TextValue $ws.Range("D39") '11.80'
$ws.Range("E39").Value = '  -3.16%  '

Set-TextValue $ws.Range("D40") '2.03'
$ws.Range("E40").Value = '  -2.78%  '

$ws.Range("E41").Value = '  +3.38%  '

Set-TextValue $ws.Range("D42") '58.45'
$ws.Range("E42").Value = '  -0.60%  '

$ws.Range("E43").Value = '  -4.61%  '

Set-TextValue $ws.Range("D44") '101.75'
$ws.Range("E44").Value = '  +4.84%  '

$ws.Range("E45").Value = '  +14.40%  '

$ws.Range("E46").Value = '  -3.15%  '

$ws.Range("E47").Value = '  -0.43%  '

$ws.Range("E48").Value = '  +9.15%  '

$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("E50").Value = '  +0.47%  '

$ws.Range("E51").Value = '  +0.86%  '

